$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 265 and 266 (pushes the existing rows
# 265-287 down to 267-289, matching the dimension change to A1:R289).
$ws.Rows.Item(265).Insert()
$ws.Rows.Item(265).Insert()

# New row 265 data
$ws.Range("A265").Value = 5
$ws.Range("B265").Value = "Macroferia Regional de Talca"
$ws.Range("C265").Value = "Maule"
$ws.Range("D265").Value = 45223
$ws.Range("E265").Value = 7
$ws.Range("F265").Value = 100112031
$ws.Range("G265").Value = "Poroto verde"
$ws.Range("H265").Value = "Sin especificar"
$ws.Range("I265").Value = "Primera"
$ws.Range("J265").Value = 100
$ws.Range("K265").Value = 25000
$ws.Range("L265").Value = 25000
$ws.Range("M265").Value = 25000
$ws.Range("N265").Value = "$/malla 25 kilos"
$ws.Range("O265").Value = "Perú"
$ws.Range("P265").Value = 1000
$ws.Range("Q265").Value = 25
$ws.Range("R265").Value = "Hortaliza"

# New row 266 data
$ws.Range("A266").Value = 5
$ws.Range("B266").Value = "Macroferia Regional de Talca"
$ws.Range("C266").Value = "Maule"
$ws.Range("D266").Value = 45223
$ws.Range("E266").Value = 7
$ws.Range("F266").Value = 100112031
$ws.Range("G266").Value = "Poroto verde"
$ws.Range("H266").Value = "Sin especificar"
$ws.Range("I266").Value = "Primera"
$ws.Range("J266").Value = 100
$ws.Range("K266").Value = 38000
$ws.Range("L266").Value = 38000
$ws.Range("M266").Value = 38000
$ws.Range("N266").Value = "$/saco 25 kilos"
$ws.Range("O266").Value = "Provincia del Elquí"
$ws.Range("P266").Value = 1520
$ws.Range("Q266").Value = 25
$ws.Range("R266").Value = "Hortaliza"
